$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 439, shifting existing rows 439:461 down to 440:462
$ws.Rows("439:439").Insert()

# Fill in the new row 439 with fresh data (matching style/format of surrounding rows)
$ws.Range("A439").Value = 10
$ws.Range("B439").Value = "Vega Modelo de Temuco"
$ws.Range("C439").Value = "La Araucanía"
$ws.Range("D439").Value = 45041
$ws.Range("D439").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E439").Value = 9
$ws.Range("F439").Value = 100112017
$ws.Range("G439").Value = "Apio"
$ws.Range("H439").Value = "Americana (o)"
$ws.Range("I439").Value = "Primera"
$ws.Range("J439").Value = 170
$ws.Range("K439").Value = 9000
$ws.Range("L439").Value = 10000
$ws.Range("M439").Value = 9471
$ws.Range("N439").Value = "$/docena de matas"
$ws.Range("O439").Value = "Provincia del Elquí"
$ws.Range("P439").Value = 1578
$ws.Range("Q439").Value = 6
$ws.Range("R439").Value = "Hortaliza"
